# LOM3238.xlsx update
# - Insert two rows after row 12 ("Docentes responsáveis:") to hold a
#   second professor entry (each on its own row, columns B/C only).
# - Fix up the text that slid into the wrong cells as a result (the
#   original file had several mismatched B/C values), and add the new
#   Portuguese translations / bibliography content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert two new rows at position 13 ---------------------------------
# This pushes the old rows 13-24 down to 15-26, which is exactly the
# row layout required by the target sheet (dimension A1:C26).
$ws.Rows("13:14").Insert()

# The inserted rows copy row 12's formatting into column A; that column
# must stay empty in the new rows (they are plain continuation rows for
# "Docentes responsáveis:"), so clear it back out.
$ws.Range("A13").Clear()
$ws.Range("A14").Clear()

# Give B13:C14 the same look (font/alignment/wrap) as the other B/C data
# cells by cloning row 15's (a real data row) formatting, then writing
# the actual values on top of it.
$ws.Range("B15:C15").Copy()
$ws.Range("B13:C14").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 2. Write the correct cell values --------------------------------------

# Objetivos: (row 10) — new Portuguese objectives text
$ws.Range("B10").Value = "Introduzir aos estudantes os princípios e a metodologia da pesquisa científica."
$ws.Range("C10").Value = "Introduzir aos estudantes os princípios e a metodologia da pesquisa científica."

# Docentes responsáveis: (rows 13-14, one professor per row)
$ws.Range("B13").Value = "519033 - Carlos Yujiro Shigue"
$ws.Range("C13").Value = "519033 - Carlos Yujiro Shigue"
$ws.Range("B14").Value = "5817692 - Katia Cristiane Gandolpho Candioto"
$ws.Range("C14").Value = "5817692 - Katia Cristiane Gandolpho Candioto"

# Programa resumido: (row 15)
$ws.Range("B15").Value = "Iniciação a um projeto de pesquisa sob orientação de um professor."
$ws.Range("C15").Value = "Iniciação a um projeto de pesquisa sob orientação de um professor."

# Programa: (row 17)
$programa = "Organização e o formalismo do desenvolvimento do trabalho científico. Técnicas de redação científica, uso de ferramentas de busca, referências bibliográficas e estruturas formais de divulgação científica. Desenvolvimento de um tema de pesquisa individual, com o formato de um trabalho de Iniciação Científica, sob a orientação de um professor ou pesquisador autorizado pela Comissão de Curso. Entrega e apresentação de monografia no final da disciplina."
$ws.Range("B17").Value = $programa
$ws.Range("C17").Value = $programa

# Método: (row 20)
$metodo = "Aulas expositivas, reuniões com professor orientador, desenvolvimento de projeto de pesquisa e elaboração de projeto de pesquisa."
$ws.Range("B20").Value = $metodo
$ws.Range("C20").Value = $metodo

# Critério: (row 21)
$criterio = "Nota de avaliação do projeto e demais documentos."
$ws.Range("B21").Value = $criterio
$ws.Range("C21").Value = $criterio

# Norma de recuperação: (row 22)
$norma = "Devido às características do curso, não será oferecida recuperação."
$ws.Range("B22").Value = $norma
$ws.Range("C22").Value = $norma

# Bibliografia: (row 23) — new bibliography text (multi-line)
$bib = "ASTI VERA, A. Metodologia da pesquisa científica. Porto Alegre: Ed. Globo, 1973.`nBARRAS, R. Os cientistas precisam escrever: guia de redação para cientistas, engenheiros e estudantes. São Paulo: TAQ/EDUSP, 1979.`nCERVO, A. L.; BERVIAN, P. A. Metodologia científica. São Paulo: Mc-Graw-Hill do Brasil, 1973.`nANDRADE, M. M. Introdução à Metodologia do Trabalho Científico São Paulo: Atlas, 2005."
$ws.Range("B23").Value = $bib
$ws.Range("C23").Value = $bib

# Note: the canonical file also narrows column A's <col> definition from
# min="1" max="2" to min="1" max="1" (column B already owns a more
# specific, later <col> entry with its real 60.71 width, so this is a
# pure metadata cleanup with no visible effect - column B's effective
# width is unchanged either way).
